$wb = $excel.ActiveWorkbook

# Use an existing header cell / data cell as a format donor so the new
# sheet's cells end up on the same style indices (s="1" header style,
# default body style) as the other sheets instead of minting brand new
# cellXfs entries for every formatting call.
$wsDeposit = $wb.Worksheets.Item("存款")
$headerDonor = $wsDeposit.Range("B1")
$indexDonor  = $wsDeposit.Range("A2")
# "2012-04-13" must stay literal text (it's the report date, not a real
# date value) - column I of 存款 already holds that exact shared string,
# so copy/paste-values from there instead of letting Range.Value auto-detect
# the string as a date serial.
$dateDonor = $wsDeposit.Range("I2")

# New worksheet "債務" (debt), placed after the last existing sheet
# ("存款") so it lands at the end of the tab strip with sheetId 5.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "債務"

# Header row (row 1), columns B..N.
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$headerDonor.Copy()
$ws.Range("B1:N1").PasteSpecial(-4122)

# Data rows 2-3, filled column-by-column (not row-by-row) so that newly
# interned shared strings land on the same indices the source workbook
# used (species, then owner, then debtor x2, then register_date x2, ...).
$ws.Range("A2").Value = 90
$ws.Range("A3").Value = 91

$ws.Range("B2").Value = "房屋貸款"
$ws.Range("B3").Value = "房屋貸款"

$ws.Range("C2").Value = "潘孟安"
$ws.Range("C3").Value = "潘孟安"

$ws.Range("D2").Value = "第一銀行恆春分行屏東縣髓鎮中正路"
$ws.Range("D3").Value = "彰化銀行賴分行屏東縣車城鄉福興村中山路"

$ws.Range("E2").Value = 2586823
$ws.Range("E3").Value = 5000000

$ws.Range("F2").Value = "98年12月24日"
$ws.Range("F3").Value = "98年07月06日"

$ws.Range("G2").Value = "房貸"
$ws.Range("G3").Value = "房貸"

$ws.Range("H2").Value = "debt"
$ws.Range("H3").Value = "debt"

$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"

$dateDonor.Copy()
$ws.Range("J2").PasteSpecial(-4163)
$dateDonor.Copy()
$ws.Range("J3").PasteSpecial(-4163)

$ws.Range("K2").Value = "潘孟安"
$ws.Range("K3").Value = "潘孟安"

$ws.Range("L2").Value = 1376
$ws.Range("L3").Value = 1376

$ws.Range("M2").Value = "tmpb07a1"
$ws.Range("M3").Value = "tmpb07a1"

$ws.Range("N2").Value = 90
$ws.Range("N3").Value = 91

# Column A on the data rows carries the same (bold/bordered/centered)
# style as the header row in every other sheet of this workbook.
$indexDonor.Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

# Restore the original active sheet so the tab-selection state of the
# pre-existing sheets is left untouched.
$wb.Worksheets.Item(1).Activate()
